$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect the Price/Volume columns from Excel's automatic number/date
# inference so values like "1.001" or "0.5100" stay text, matching the
# inlineStr cells already in the sheet.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.782.60"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "1.624.97"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "214.91"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "0.5100"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.2557"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "0.06316"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").Value = "19.35"
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("D11").Value = "0.07764"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "4.220"
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("D13").Value = "1.624.37"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "1.845.11"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").Value = "0.5513"
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").Value = "63.42"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "0.0₅7491"
$ws.Range("E17").Value = "  -2.46%  "
$ws.Range("D18").Value = "25.785.77"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "193.96"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "4.401"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "9.749"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").Value = "5.995"
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "1.871"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "141.37"
$ws.Range("D27").Value = "0.1242"
$ws.Range("E27").Value = "  +3.97%  "
$ws.Range("D28").Value = "15.51"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").Value = "6.695"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("D30").Value = "1.237"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "0.04848"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").Value = "3.230"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").Value = "3.151"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D35").Value = "2.366"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "0.8925"
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("D37").Value = "2.537"
$ws.Range("E37").Value = "  -1.81%  "
$ws.Range("D38").Value = "0.5489"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("D39").Value = "1.113.19"
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").Value = "0.01544"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("D43").Value = "0.7954"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").Value = "96.96"
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("D45").Value = "1.770.50"
$ws.Range("E46").Value = "  -7.22%  "
$ws.Range("D47").Value = "0.4421"
$ws.Range("E47").Value = "  -2.35%  "
$ws.Range("D48").Value = "0.9958"
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("D49").Value = "54.46"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").Value = "0.05126"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").Value = "7.542"
$ws.Range("E51").Value = "  +3.23%  "

# Reapply the default style so cells don't carry a leftover text
# number-format override (matches original formatting).
$ws.Range("D2:E51").Style = "Normal"
